# GuestConfig via DSC - "Apply and Set" deck
# Slide 3, shape "Rectangle 3": retitle the banner text and give the
# text frame its (now explicit) default insets.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(3)
$sh = $s.Shapes.Item("Rectangle 3")

# --- body insets: make the (already-default) margins explicit ------------
$tf = $sh.TextFrame
$tf.MarginLeft   = 7.2   # 91440 EMU
$tf.MarginTop    = 3.6   # 45720 EMU
$tf.MarginRight  = 7.2   # 91440 EMU
$tf.MarginBottom = 3.6   # 45720 EMU

# --- retitle: split into 3 runs (same rPr) around "CheckForService" -------
$tr = $sh.TextFrame.TextRange

$part1 = "Build Authoring VM | Author DSC "
$part2 = "CheckForService"
$enDash = [char]0x2013
$part3 = " " + $enDash + " Set WDAV Passive"

$tr.Text = $part1 + $part2 + $part3

# Re-assigning the same text over each sub-range forces PowerPoint to
# split the run at that boundary without touching any formatting.
$sub2 = $tr.Characters($part1.Length + 1, $part2.Length)
$sub2.Text = $part2

$sub3 = $tr.Characters($part1.Length + $part2.Length + 1, $part3.Length)
$sub3.Text = $part3
